$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear rows 6 and 7 entirely (they become unused/gap rows)
$ws.Rows(6).Clear()
$ws.Rows(7).Clear()

# Row 2
$ws.Range('A2').Value = '904845630'
$ws.Range('B2').Value = '20171025'
$ws.Range('C2').Value = '3013996080'
$ws.Range('D2').Value = '06B'
$ws.Range('E2').Value = '14446757'
$ws.Range('F2').Value = 'WASHINGTON DUKE     '
$ws.Range('G2').Value = 'USEE000001'
$ws.Range('H2').Value = 'TA5TVWRLR      '
$ws.Range('I2').Value = 'Pro Wheeled Roller Bag Black  '
$ws.Range('J2').Value = '1'
$ws.Range('K2').Value = 'EA '
$ws.Range('L2').Value = '180.00'
$ws.Range('M2').Value = '0.00'
$ws.Range('N2').Value = '0.00'
$ws.Range('O2').Value = '100.00'
$ws.Range('P2').Value = '0.00'
$ws.Range('Q2').Value = '0.00'
$ws.Range('R2').Value = '0.00'
$ws.Range('S2').Value = '1'
$ws.Range('T2').Value = '0.00'
$ws.Range('U2').Value = '0.00'
$ws.Range('V2').Value = ' 0.00'

# Row 3
$ws.Range('A3').Value = '904845629'
$ws.Range('B3').Value = '20171025'
$ws.Range('C3').Value = '3013996079'
$ws.Range('D3').Value = '06B'
$ws.Range('E3').Value = '14446756'
$ws.Range('F3').Value = '                    '
$ws.Range('G3').Value = 'USEM000044'
$ws.Range('H3').Value = 'TA5TVWRLR      '
$ws.Range('I3').Value = 'Pro Wheeled Roller Bag Black  '
$ws.Range('J3').Value = '1'
$ws.Range('K3').Value = 'EA '
$ws.Range('L3').Value = '180.00'
$ws.Range('M3').Value = '0.00'
$ws.Range('N3').Value = '0.00'
$ws.Range('O3').Value = '100.00'
$ws.Range('P3').Value = '0.00'
$ws.Range('Q3').Value = '0.00'
$ws.Range('R3').Value = '0.00'
$ws.Range('S3').Value = '1'
$ws.Range('T3').Value = '0.00'
$ws.Range('U3').Value = '0.00'
$ws.Range('V3').Value = ' 0.00'

# Row 4
$ws.Range('A4').Value = '904845652'
$ws.Range('B4').Value = '20171025'
$ws.Range('C4').Value = '3013996239'
$ws.Range('D4').Value = '06B'
$ws.Range('E4').Value = '14446940'
$ws.Range('F4').Value = 'CART MITTS          '
$ws.Range('G4').Value = 'US00059139'
$ws.Range('H4').Value = 'TA7WEACM-0     '
$ws.Range('I4').Value = 'Titleist Cart Mitts           '
$ws.Range('J4').Value = '6'
$ws.Range('K4').Value = 'EA '
$ws.Range('L4').Value = '22.50'
$ws.Range('M4').Value = '0.00'
$ws.Range('N4').Value = '0.00'
$ws.Range('O4').Value = '0.00'
$ws.Range('P4').Value = '0.00'
$ws.Range('Q4').Value = '22.50'
$ws.Range('R4').Value = '135.00'
$ws.Range('S4').Value = '6'
$ws.Range('T4').Value = '135.00'
$ws.Range('U4').Value = '9.98'
$ws.Range('V4').Value = ' 135.00'

# Row 5
$ws.Range('A5').Value = '904845651'
$ws.Range('B5').Value = '20171025'
$ws.Range('C5').Value = '3013996237'
$ws.Range('D5').Value = '06B'
$ws.Range('E5').Value = '14446939'
$ws.Range('F5').Value = 'WINTER CAPS         '
$ws.Range('G5').Value = 'US00059139'
$ws.Range('H5').Value = 'TH7WEALB-P06   '
$ws.Range('I5').Value = 'Lifestyle Beanie Legacy Asst  '
$ws.Range('J5').Value = '6'
$ws.Range('K5').Value = 'EA '
$ws.Range('L5').Value = '12.50'
$ws.Range('M5').Value = '0.00'
$ws.Range('N5').Value = '0.00'
$ws.Range('O5').Value = '0.00'
$ws.Range('P5').Value = '0.00'
$ws.Range('Q5').Value = '12.50'
$ws.Range('R5').Value = '75.00'
$ws.Range('S5').Value = '6'
$ws.Range('T5').Value = '75.00'
$ws.Range('U5').Value = '8.39'
$ws.Range('V5').Value = ' 75.00'

# Row 8
$ws.Range('A8').Value = '904845552'
$ws.Range('B8').Value = '20171025'
$ws.Range('C8').Value = '3013953706'
$ws.Range('D8').Value = '06B'
$ws.Range('E8').Value = '14258881'
$ws.Range('F8').Value = '                    '
$ws.Range('G8').Value = 'US00026239'
$ws.Range('H8').Value = 'TH7WEAWHP-P06  '
$ws.Range('I8').Value = 'Pom Pom Winter Hat Asst.      '
$ws.Range('J8').Value = '6'
$ws.Range('K8').Value = 'EA '
$ws.Range('L8').Value = '15.00'
$ws.Range('M8').Value = '0.00'
$ws.Range('N8').Value = '0.00'
$ws.Range('O8').Value = '0.00'
$ws.Range('P8').Value = '0.00'
$ws.Range('Q8').Value = '15.00'
$ws.Range('R8').Value = '90.00'
$ws.Range('S8').Value = '6'
$ws.Range('T8').Value = '90.00'
$ws.Range('U8').Value = '6.39'
$ws.Range('V8').Value = ' 90.00'

# Row 9
$ws.Range('A9').Value = '904845557'
$ws.Range('B9').Value = '20171025'
$ws.Range('C9').Value = '3013965607'
$ws.Range('D9').Value = '06B'
$ws.Range('E9').Value = '14320931'
$ws.Range('F9').Value = 'POM/POMS            '
$ws.Range('G9').Value = 'US00033077'
$ws.Range('H9').Value = 'TH7WEAWHP-P06  '
$ws.Range('I9').Value = 'Pom Pom Winter Hat Asst.      '
$ws.Range('J9').Value = '12'
$ws.Range('K9').Value = 'EA '
$ws.Range('L9').Value = '15.00'
$ws.Range('M9').Value = '0.00'
$ws.Range('N9').Value = '0.00'
$ws.Range('O9').Value = '0.00'
$ws.Range('P9').Value = '0.00'
$ws.Range('Q9').Value = '15.00'
$ws.Range('R9').Value = '180.00'
$ws.Range('S9').Value = '12'
$ws.Range('T9').Value = '180.00'
$ws.Range('U9').Value = '8.78'
$ws.Range('V9').Value = ' 180.00'

# Row 10
$ws.Range('A10').Value = '904845587'
$ws.Range('B10').Value = '20171025'
$ws.Range('C10').Value = '3013996644'
$ws.Range('D10').Value = '06Y'
$ws.Range('E10').Value = '14447662'
$ws.Range('F10').Value = 'CUSTOM              '
$ws.Range('G10').Value = 'US00002181'
$ws.Range('H10').Value = 'TA5ACMFTWLC    '
$ws.Range('I10').Value = 'Waffle Microfiber Towel CST   '
$ws.Range('J10').Value = '1'
$ws.Range('K10').Value = 'EA '
$ws.Range('L10').Value = '17.00'
$ws.Range('M10').Value = '0.00'
$ws.Range('N10').Value = '0.00'
$ws.Range('O10').Value = '0.00'
$ws.Range('P10').Value = '0.00'
$ws.Range('Q10').Value = '17.00'
$ws.Range('R10').Value = '17.00'
$ws.Range('S10').Value = '1'
$ws.Range('T10').Value = '17.00'
$ws.Range('U10').Value = '9.07'
$ws.Range('V10').Value = ' 17.00'

# Row 11
$ws.Range('A11').Value = '904845551'
$ws.Range('B11').Value = '20171025'
$ws.Range('C11').Value = '3013953151'
$ws.Range('D11').Value = '06B'
$ws.Range('E11').Value = '14257783'
$ws.Range('F11').Value = '                    '
$ws.Range('G11').Value = 'US00031646'
$ws.Range('H11').Value = 'TA7WEACM-0     '
$ws.Range('I11').Value = 'Titleist Cart Mitts           '
$ws.Range('J11').Value = '6'
$ws.Range('K11').Value = 'EA '
$ws.Range('L11').Value = '22.50'
$ws.Range('M11').Value = '0.00'
$ws.Range('N11').Value = '0.00'
$ws.Range('O11').Value = '0.00'
$ws.Range('P11').Value = '0.00'
$ws.Range('Q11').Value = '22.50'
$ws.Range('R11').Value = '135.00'
$ws.Range('S11').Value = '6'
$ws.Range('T11').Value = '135.00'
$ws.Range('U11').Value = '9.11'
$ws.Range('V11').Value = ' 135.00'

# Row 12
$ws.Range('A12').Value = '904845583'
$ws.Range('B12').Value = '20171025'
$ws.Range('C12').Value = '3013995773'
$ws.Range('D12').Value = '06B'
$ws.Range('E12').Value = '14446306'
$ws.Range('F12').Value = 'S                   '
$ws.Range('G12').Value = 'US00060065'
$ws.Range('H12').Value = 'TH7VTP-P12     '
$ws.Range('I12').Value = 'Tour Perf Visor Legacy Asst   '
$ws.Range('J12').Value = '12'
$ws.Range('K12').Value = 'EA '
$ws.Range('L12').Value = '12.00'
$ws.Range('M12').Value = '0.00'
$ws.Range('N12').Value = '0.00'
$ws.Range('O12').Value = '0.00'
$ws.Range('P12').Value = '0.00'
$ws.Range('Q12').Value = '8.00'
$ws.Range('R12').Value = '96.00'
$ws.Range('S12').Value = '12'
$ws.Range('T12').Value = '96.00'
$ws.Range('U12').Value = '8.78'
$ws.Range('V12').Value = ' 96.00'

# Row 13
$ws.Range('A13').Value = '904845581'
$ws.Range('B13').Value = '20171025'
$ws.Range('C13').Value = '3013995485'
$ws.Range('D13').Value = '06B'
$ws.Range('E13').Value = '14445845'
$ws.Range('F13').Value = 'THOMAS              '
$ws.Range('G13').Value = 'US00002681'
$ws.Range('H13').Value = 'TH7ASC-P06     '
$ws.Range('I13').Value = 'Tour Snapback Asst.           '
$ws.Range('J13').Value = '6'
$ws.Range('K13').Value = 'EA '
$ws.Range('L13').Value = '15.00'
$ws.Range('M13').Value = '0.00'
$ws.Range('N13').Value = '0.00'
$ws.Range('O13').Value = '0.00'
$ws.Range('P13').Value = '0.00'
$ws.Range('Q13').Value = '15.00'
$ws.Range('R13').Value = '90.00'
$ws.Range('S13').Value = '6'
$ws.Range('T13').Value = '90.00'
$ws.Range('U13').Value = '8.39'
$ws.Range('V13').Value = ' 90.00'

# Row 14
$ws.Range('A14').Value = '904845578'
$ws.Range('B14').Value = '20171025'
$ws.Range('C14').Value = '3013995270'
$ws.Range('D14').Value = '06B'
$ws.Range('E14').Value = '14445144'
$ws.Range('F14').Value = 'WLA001716-27        '
$ws.Range('G14').Value = 'US00057861'
$ws.Range('H14').Value = 'TA1ACFVP-0     '
$ws.Range('I14').Value = 'Fleece Valuables Pouch        '
$ws.Range('J14').Value = '4'
$ws.Range('K14').Value = 'EA '
$ws.Range('L14').Value = '6.00'
$ws.Range('M14').Value = '0.00'
$ws.Range('N14').Value = '2.00'
$ws.Range('O14').Value = '0.00'
$ws.Range('P14').Value = '0.00'
$ws.Range('Q14').Value = '5.88'
$ws.Range('R14').Value = '23.52'
$ws.Range('S14').Value = '4'
$ws.Range('T14').Value = '23.52'
$ws.Range('U14').Value = '6.00'
$ws.Range('V14').Value = ' 23.52'

# Update selection to match target (entire row 8 selected)
$ws.Rows(8).Select()
